# Updates cryptos list: Price (column D) and Volume(1h) (column E) cells
# for the rows whose market data changed, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. NumberFormat is forced to
# '@' (Text) before each write so numeric-looking strings such as
# "0.998" or "4.00" are stored as text (matching the original
# inlineStr cells) instead of being coerced into numbers.
$updates = [ordered]@{
    'D2' = '60.549.48'
    'E2' = '  +0.12%  '
    'D3' = '2.588.64'
    'E3' = '  +0.21%  '
    'E4' = '  -0.27%  '
    'D5' = '508.12'
    'E5' = '  -0.90%  '
    'D6' = '154.97'
    'E6' = '  -2.37%  '
    'D7' = '0.998'
    'E7' = '  +0.76%  '
    'D8' = '0.581'
    'E8' = '  -4.53%  '
    'D9' = '2.597.85'
    'E9' = '  -1.76%  '
    'E10' = '  +6.51%  '
    'E11' = '  -1.10%  '
    'D12' = '0.347'
    'E12' = '  +0.02%  '
    'E13' = '  +1.23%  '
    'D14' = '3.044.65'
    'E14' = '  +0.35%  '
    'D15' = '60.606.62'
    'E15' = '  +0.70%  '
    'D16' = '21.68'
    'E16' = '  -2.25%  '
    'D17' = '0.0000140'
    'E17' = '  -0.36%  '
    'D18' = '2.592.18'
    'E18' = '  -0.95%  '
    'E19' = '  -1.09%  '
    'D20' = '346.98'
    'E20' = '  +0.84%  '
    'D21' = '10.49'
    'E21' = '  -0.83%  '
    'D22' = '6.12'
    'E22' = '  -0.73%  '
    'D23' = '0.997'
    'E23' = '  +0.01%  '
    'D24' = '60.05'
    'E24' = '  -0.96%  '
    'E25' = '  -0.81%  '
    'E26' = '  -1.15%  '
    'D27' = '2.699.75'
    'E27' = '  +0.26%  '
    'D28' = '0.999'
    'E28' = '  +0.92%  '
    'D29' = '0.0₃0847'
    'E29' = '  -2.07%  '
    'D30' = '7.41'
    'E30' = '  -1.70%  '
    'E31' = '  +0.42%  '
    'D32' = '19.41'
    'E32' = '  -1.38%  '
    'D33' = '153.21'
    'E33' = '  -2.65%  '
    'E34' = '  -1.74%  '
    'D35' = '5.72'
    'E35' = '  +2.39%  '
    'D36' = '4.00'
    'E36' = '  +0.44%  '
    'E37' = '  -2.09%  '
    'E38' = '  +7.00%  '
    'E39' = '  -0.13%  '
    'D40' = '0.849'
    'E40' = '  -2.75%  '
    'D41' = '36.17'
    'E41' = '  +1.69%  '
    'D43' = '296.73'
    'E43' = '  -2.91%  '
    'D44' = '0.624'
    'E44' = '  -2.56%  '
    'D45' = '0.0993'
    'E45' = '  -2.18%  '
    'D46' = '0.0559'
    'E46' = '  -3.08%  '
    'D47' = '0.998'
    'E47' = '  +0.89%  '
    'D48' = '19.87'
    'E48' = '  -0.22%  '
    'E49' = '  -2.70%  '
    'D50' = '0.0233'
    'E50' = '  -2.08%  '
    'D51' = '10.30'
    'E51' = '  +0.15%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
